$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'sliding pads for exercise',
    'compression capri men',
    'basketball pants for boys',
    'padded baseball pants',
    'cycling capri pants',
    'sliding pants baseball mens',
    'below the knee shorts for men',
    'mens long basketball shorts below knee',
    'youth small knee pads basketball',
    'basketball leg pads',
    'short baseball pants',
    'basketball bump',
    'knee pads running',
    'football compression shorts with pads',
    'knee high baseball pants mens',
    'knee padded compression',
    'athletic pads',
    'mens protection pads',
    'youth knee pads wrestling',
    'men softball pants',
    'baseball sliding pants',
    'boys sliding pants',
    'youth boys basketball pants',
    'lacrosse sweat pants',
    'raceface knee pads',
    'capri compression pants men',
    'snowboarding compression pants',
    'men knee pad pants',
    'youth basketball tights for boys',
    'knee pads for basketball youth',
    'mtn bike knee pads',
    'baseball padded sliding shorts',
    'leggings with baseballs',
    'youth knee pad wrestling',
    '28 basketball',
    'tights with knee',
    'compression capris',
    'padded tights for football',
    'baseball tights for boys',
    'baseball sliding shorts men',
    'knee pads for basketball youth boys',
    'youth leggings boys basketball',
    'basketballs leggings',
    'softball items',
    'football leggings for men',
    'basketball compression pants youth',
    'compression pants men black',
    'mens work pants with knee pads',
    'youth tights',
    'mens down pants',
    'compression basketball pants youth',
    'basketball compression pants women',
    'softball slider',
    'work knee pads under pants',
    'calf compression pants',
    'leg pads basketball',
    'compression pants men football',
    'softball shorts men',
    'mens softball shorts',
    'mens capri shorts below knee',
    'compression tights youth',
    'adidas knee pads',
    'athletic leggings mens',
    'knee compression shorts',
    'compression knee pads pair',
    'indoor volleyball knee pads',
    'youth basketball pants boys',
    'softball sliding shorts girls padded',
    'lacrosse compression shorts padded',
    'taken leggings',
    'mens basketball tights',
    'black compression pants men',
    'softball sliding shorts',
    'girls basketball knee pads youth',
    'kneepads basketball',
    'knee protector for construction',
    'knee pad sleeve basketball',
    'youth basketball leggings',
    'medium compression pants',
    'football knee pads for men',
    'mens compression leggings',
    'knee pad for yoga',
    'boys knee pads basketball',
    'baseball sliding',
    'knee pads for biking men',
    'knee pads girls basketball',
    'black football pants',
    'lacrosse pants',
    'team work softball pants',
    'long basketball shorts for men below knee',
    'basketball hex pads',
    'compression shorts padded basketball',
    'big boys tights',
    'basketball tights',
    'hockey knee pads adult',
    'padded compression pants football',
    'hockey hip pads adult',
    'soccer goalkeeper pads',
    'athletic mens leggings',
    'knee pads for basketball youth'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
